$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values formatted as plain text (e.g. "42.474.63"),
# including values that look numeric (e.g. "72.20"). Force text format so
# Excel does not silently coerce them to numbers and drop trailing zeros.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '42.531.77'
$ws.Range('E2').Value = '  -2.96%  '
$ws.Range('D3').Value = '2.256.70'
$ws.Range('E3').Value = '  -3.99%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '233.52'
$ws.Range('E5').Value = '  -2.98%  '
$ws.Range('D6').Value = '0.639'
$ws.Range('E6').Value = '  -3.84%  '
$ws.Range('D7').Value = '72.20'
$ws.Range('E7').Value = '  -1.29%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').Value = '0.569'
$ws.Range('E9').Value = '  -4.86%  '
$ws.Range('D10').Value = '0.0998'
$ws.Range('E10').Value = '  -0.55%  '
$ws.Range('D11').Value = '59.06'
$ws.Range('E11').Value = '  -2.95%  '
$ws.Range('D12').Value = '36.92'
$ws.Range('E12').Value = '  +12.19%  '
$ws.Range('E13').Value = '  -2.25%  '
$ws.Range('D14').Value = '6.91'
$ws.Range('E14').Value = '  -5.63%  '
$ws.Range('D15').Value = '2.583.45'
$ws.Range('E15').Value = '  -4.34%  '
$ws.Range('D16').Value = '15.17'
$ws.Range('E16').Value = '  -7.62%  '
$ws.Range('D17').Value = '0.881'
$ws.Range('E17').Value = '  -2.78%  '
$ws.Range('D18').Value = '2.244.28'
$ws.Range('E18').Value = '  -4.47%  '
$ws.Range('D19').Value = '42.409.20'
$ws.Range('E19').Value = '  -3.07%  '
$ws.Range('D20').Value = '0.0₃0995'
$ws.Range('E20').Value = '  -2.40%  '
$ws.Range('D21').Value = '74.03'
$ws.Range('E21').Value = '  -4.01%  '
$ws.Range('D22').Value = '6.27'
$ws.Range('E22').Value = '  -6.23%  '
$ws.Range('D23').Value = '238.88'
$ws.Range('E23').Value = '  -6.83%  '
$ws.Range('D24').Value = '1.96'
$ws.Range('E24').Value = '  +2.52%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').Value = '3.69'
$ws.Range('E26').Value = '  -1.03%  '
$ws.Range('D27').Value = '2.38'
$ws.Range('E27').Value = '  -4.94%  '
$ws.Range('D28').Value = '10.17'
$ws.Range('E28').Value = '  -3.83%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.17'
$ws.Range('E29').Value = '  -8.26%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '168.49'
$ws.Range('E30').Value = '  -4.99%  '
$ws.Range('E31').Value = '  -7.89%  '
$ws.Range('D32').Value = '0.122'
$ws.Range('E32').Value = '  -5.35%  '
$ws.Range('D33').Value = '0.128'
$ws.Range('E33').Value = '  -5.73%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.0732'
$ws.Range('E34').Value = '  -3.06%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '5.42'
$ws.Range('E35').Value = '  -0.80%  '
$ws.Range('D36').Value = '4.85'
$ws.Range('E36').Value = '  -6.48%  '
$ws.Range('D37').Value = '3.68'
$ws.Range('E37').Value = '  -3.76%  '
$ws.Range('D38').Value = '22.70'
$ws.Range('E38').Value = '  +18.96%  '
$ws.Range('D39').Value = '6.17'
$ws.Range('E39').Value = '  -2.72%  '
$ws.Range('D40').Value = '2.27'
$ws.Range('E40').Value = '  -4.52%  '
$ws.Range('E41').Value = '  -3.61%  '
$ws.Range('D42').Value = '66.45'
$ws.Range('E42').Value = '  -3.02%  '
$ws.Range('D43').Value = '5.20'
$ws.Range('E43').Value = '  +4.76%  '
$ws.Range('D44').Value = '8.97'
$ws.Range('E44').Value = '  -1.86%  '
$ws.Range('E45').Value = '  -7.12%  '
$ws.Range('D46').Value = '0.192'
$ws.Range('E46').Value = '  -4.83%  '
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('D48').Value = '4.46'
$ws.Range('E48').Value = '  +6.04%  '
$ws.Range('D49').Value = '2.41'
$ws.Range('E49').Value = '  -3.51%  '
$ws.Range('B50').Value = 'Celestia'
$ws.Range('C50').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D50').Value = '10.32'
$ws.Range('E50').Value = '  +9.32%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').Value = '1.20'
$ws.Range('E51').Value = '  -4.10%  '
